$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 82.083336
$ws.Range("I11").Value = 82.083336
$ws.Range("K11").Value = 82.083336
$ws.Range("M11").Value = 57.916664
$ws.Range("H17").Value = 2641.5
$ws.Range("I17").Value = 2487.5
$ws.Range("J17").Value = 2949.5
$ws.Range("K17").Value = 7462.5
$ws.Range("L17").Value = 8848.5
$ws.Range("M17").Value = -7294.5
$ws.Range("N17").Value = -9184.5
$ws.Range("H32").Value = 2683
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2683
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 2683
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -3335
$ws.Range("H33").Value = 181.22223
$ws.Range("I33").Value = 197.35715
$ws.Range("K33").Value = 197.35715
$ws.Range("M33").Value = 31.64285000000001
$ws.Range("H87").Value = 54214.285
$ws.Range("J87").Value = 54214.285
$ws.Range("L87").Value = 54214.285
$ws.Range("N87").Value = -56710.285
$ws.Range("H90").Value = 54214.285
$ws.Range("J90").Value = 54214.285
$ws.Range("L90").Value = 162642.855
$ws.Range("N90").Value = -175122.855
$ws.Range("H132").Value = 3607.353
$ws.Range("I132").Value = 3607.353
$ws.Range("K132").Value = 10822.059
$ws.Range("M132").Value = -8292.059000000001
$ws.Range("H141").Value = 816.2
$ws.Range("I141").Value = 816.2
$ws.Range("K141").Value = 2448.6
$ws.Range("M141").Value = 2731.4

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 7650
$ws.Range("I63").Value = 7475
$ws.Range("K63").Value = 7475
$ws.Range("M63").Value = -6789
$ws.Range("H66").Value = 7650
$ws.Range("I66").Value = 7475
$ws.Range("K66").Value = 37375
$ws.Range("M66").Value = -33943
$ws.Range("H74").Value = 4827.9
$ws.Range("I74").Value = 4626
$ws.Range("J74").Value = 4962.5
$ws.Range("K74").Value = 4626
$ws.Range("L74").Value = 4962.5
$ws.Range("M74").Value = -3752
$ws.Range("N74").Value = -6710.5
$ws.Range("H77").Value = 4827.9
$ws.Range("I77").Value = 4626
$ws.Range("J77").Value = 4962.5
$ws.Range("K77").Value = 23130
$ws.Range("L77").Value = 24812.5
$ws.Range("M77").Value = -18762
$ws.Range("N77").Value = -33548.5
$ws.Range("H92").Value = 61385.75
$ws.Range("I92").Value = 90000
$ws.Range("J92").Value = 51847.668
$ws.Range("K92").Value = 90000
$ws.Range("L92").Value = 51847.668
$ws.Range("M92").Value = -87504
$ws.Range("N92").Value = -56839.668
$ws.Range("H97").Value = 3336.25
$ws.Range("I97").Value = 3384.2856
$ws.Range("K97").Value = 3384.2856
$ws.Range("M97").Value = -2888.2856
$ws.Range("H102").Value = 1999.5
$ws.Range("I102").Value = 1999
$ws.Range("K102").Value = 1999
$ws.Range("M102").Value = -377
$ws.Range("H132").Value = 8016.4165
$ws.Range("I132").Value = 6299.6665
$ws.Range("K132").Value = 18898.9995
$ws.Range("M132").Value = -16368.9995
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 894
$ws.Range("I122").Value = 894
$ws.Range("K122").Value = 2682
$ws.Range("M122").Value = -232
$ws.Range("H132").Value = 5265.375
$ws.Range("J132").Value = 13000
$ws.Range("L132").Value = 39000
$ws.Range("N132").Value = -44060

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 100001450
$ws.Range("I4").Value = 1610.2222
$ws.Range("K4").Value = 4830.6666
$ws.Range("M4").Value = -4718.6666
$ws.Range("H13").Value = 2239.6
$ws.Range("I13").Value = 2849.3333
$ws.Range("K13").Value = 8547.999899999999
$ws.Range("M13").Value = -8379.999899999999
$ws.Range("H68").Value = 1481.8
$ws.Range("J68").Value = 1455
$ws.Range("L68").Value = 4365
$ws.Range("N68").Value = -5987
$ws.Range("H71").Value = 1481.8
$ws.Range("J71").Value = 1455
$ws.Range("L71").Value = 13095
$ws.Range("N71").Value = -21207
$ws.Range("H109").Value = 3299
$ws.Range("I109").Value = 3299
$ws.Range("K109").Value = 9897
$ws.Range("M109").Value = -8857

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 3250
$ws.Range("J25").Value = 5000
$ws.Range("L25").Value = 5000
$ws.Range("N25").Value = -6058

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7598.6
$ws.Range("I16").Value = 7598.6
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 7598.6
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -7428.6
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 5964.95
$ws.Range("I22").Value = 5956.1875
$ws.Range("K22").Value = 5956.1875
$ws.Range("M22").Value = -5661.1875
$ws.Range("H27").Value = 5964.95
$ws.Range("I27").Value = 5956.1875
$ws.Range("K27").Value = 5956.1875
$ws.Range("M27").Value = -5849.1875
$ws.Range("H68").Value = 6499
$ws.Range("I68").Value = 6499
$ws.Range("K68").Value = 6499
$ws.Range("M68").Value = -5750
$ws.Range("H71").Value = 6499
$ws.Range("I71").Value = 6499
$ws.Range("K71").Value = 32495
$ws.Range("M71").Value = -28751
$ws.Range("H132").Value = 15500
$ws.Range("I132").Value = 6500
$ws.Range("J132").Value = 24500
$ws.Range("K132").Value = 19500
$ws.Range("L132").Value = 73500
$ws.Range("M132").Value = -16970
$ws.Range("N132").Value = -78560
$ws.Range("H136").Value = 17599
$ws.Range("I136").Value = 12998.333
$ws.Range("K136").Value = 38994.999
$ws.Range("M136").Value = -36444.999

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2065.6667
$ws.Range("I113").Value = 2065.6667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 6197.000100000001
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -4027.000100000001
$ws.Range("N113").ClearContents()
$ws.Range("H136").Value = 7388.1377
$ws.Range("I136").Value = 6076.478
$ws.Range("K136").Value = 18229.434
$ws.Range("M136").Value = -15679.434
